# Rename the only worksheet from "Arkusz1" to "Rodzaj eksplant na efekt"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$oldName = $ws.Name
$newName = "Rodzaj eksplant na efekt"
$ws.Name = $newName

# Update the embedded charts' series formulas so they point at the
# renamed sheet instead of the old "Arkusz1" name.
$chartObjects = $ws.ChartObjects()
for ($i = 1; $i -le $chartObjects.Count; $i++) {
    $chartObject = $chartObjects.Item($i)
    $chart = $chartObject.Chart
    $series = $chart.SeriesCollection()
    for ($j = 1; $j -le $series.Count; $j++) {
        $s = $series.Item($j)
        $s.Formula = $s.Formula.Replace($oldName + "!", "'" + $newName + "'!")
    }
}

# Move the active selection on the worksheet from G20 to C21
$ws.Range("C21").Select()
